# Auto-generated edit.ps1
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Sheet "总计": insert a new row 2 for "2022-Q3", shifting the
#    existing quarterly rows down by one.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 8; $r -ge 2; $r--) {
    $dstRow = $r + 1
    $summary.Rows("$r`:$r").Copy()
    $summary.Rows("$dstRow`:$dstRow").PasteSpecial(-4104)
    # PasteSpecial onto a brand-new row does not carry the cell
    # style for column A, so re-apply it explicitly.
    $summary.Range("A$r").Copy()
    $summary.Range("A$dstRow").PasteSpecial(-4122)
    # Column A is a 0-based row sequence number, not a copied
    # value - fix it up to match its new row position.
    $summary.Range("A$dstRow").Value = $dstRow - 2
}

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = '2022-Q3'
$summary.Range("C2").Value = 36
$summary.Range("D2").Value = 3.08

# ---------------------------------------------------------------
# 2) New sheet "2022-Q3": holdings detail, inserted right after
#    "总计" (i.e. before the old first quarter sheet).
# ---------------------------------------------------------------
$firstQtrSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($firstQtrSheet)
$q3.Name = "2022-Q3"

# Header row (copy the header style used on "总计"!B1, s=2)
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("B1").Value = '基金代码'
$q3.Range("C1").Value = '基金名称'
$q3.Range("D1").Value = '基金规模'
$q3.Range("E1").Value = '股票总仓位'
$q3.Range("F1").Value = '仓位占比'
$q3.Range("G1").Value = '持有市值(亿元)'
$q3.Range("H1").Value = '仓位排名'

# Column A (sequence number) uses the same bordered/bold style as
# on the "总计" sheet (s=2); apply it once to the whole column range
# then fill in the values below.
$summary.Range("A2").Copy()
$q3.Range("A2:A37").PasteSpecial(-4122)

# Columns B, D:G hold numeric-looking values that must stay TEXT
# (fund codes with leading zeros, and numbers-as-text as in the
# source file). Pre-format the ranges as text so assigning the
# numeric-looking strings does not get reinterpreted as a number.
$q3.Range("B2:B37").NumberFormat = "@"
$q3.Range("D2:G37").NumberFormat = "@"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = '000478'
$q3.Range("C2").Value = '建信中证500指数增强A'
$q3.Range("D2").Value = '45.95'
$q3.Range("E2").Value = '82.53'
$q3.Range("F2").Value = '1.12'
$q3.Range("G2").Value = '0.5146'
$q3.Range("H2").Value = 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = '000780'
$q3.Range("C3").Value = '鹏华医疗保健股票'
$q3.Range("D3").Value = '6.60'
$q3.Range("E3").Value = '81.50'
$q3.Range("F3").Value = '7.50'
$q3.Range("G3").Value = '0.4950'
$q3.Range("H3").Value = 1
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = '012284'
$q3.Range("C4").Value = '光大保德信健康优加混合'
$q3.Range("D4").Value = '11.36'
$q3.Range("E4").Value = '89.34'
$q3.Range("F4").Value = '3.24'
$q3.Range("G4").Value = '0.3681'
$q3.Range("H4").Value = 10
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = '005505'
$q3.Range("C5").Value = '前海开源中药研究精选股票A'
$q3.Range("D5").Value = '4.10'
$q3.Range("E5").Value = '83.96'
$q3.Range("F5").Value = '6.85'
$q3.Range("G5").Value = '0.2808'
$q3.Range("H5").Value = 9
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = '005506'
$q3.Range("C6").Value = '前海开源中药研究精选股票C'
$q3.Range("D6").Value = '3.58'
$q3.Range("E6").Value = '83.96'
$q3.Range("F6").Value = '6.85'
$q3.Range("G6").Value = '0.2452'
$q3.Range("H6").Value = 9
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = '002919'
$q3.Range("C7").Value = '东吴智慧医疗量化策略灵活配置混合A'
$q3.Range("D7").Value = '3.83'
$q3.Range("E7").Value = '90.19'
$q3.Range("F7").Value = '6.28'
$q3.Range("G7").Value = '0.2405'
$q3.Range("H7").Value = 2
$q3.Range("A8").Value = 6
$q3.Range("B8").Value = '011948'
$q3.Range("C8").Value = '东吴智慧医疗量化策略灵活配置混合C'
$q3.Range("D8").Value = '2.23'
$q3.Range("E8").Value = '90.19'
$q3.Range("F8").Value = '6.28'
$q3.Range("G8").Value = '0.1400'
$q3.Range("H8").Value = 2
$q3.Range("A9").Value = 7
$q3.Range("B9").Value = '163503'
$q3.Range("C9").Value = '天治核心成长混合（LOF）'
$q3.Range("D9").Value = '3.36'
$q3.Range("E9").Value = '93.81'
$q3.Range("F9").Value = '2.98'
$q3.Range("G9").Value = '0.1001'
$q3.Range("H9").Value = 10
$q3.Range("A10").Value = 8
$q3.Range("B10").Value = '002938'
$q3.Range("C10").Value = '中银证券健康产业灵活配置混合'
$q3.Range("D10").Value = '1.98'
$q3.Range("E10").Value = '92.72'
$q3.Range("F10").Value = '4.55'
$q3.Range("G10").Value = '0.0901'
$q3.Range("H10").Value = 7
$q3.Range("A11").Value = 9
$q3.Range("B11").Value = '013940'
$q3.Range("C11").Value = '东吴医疗服务股票A'
$q3.Range("D11").Value = '1.10'
$q3.Range("E11").Value = '91.62'
$q3.Range("F11").Value = '6.55'
$q3.Range("G11").Value = '0.0720'
$q3.Range("H11").Value = 2
$q3.Range("A12").Value = 10
$q3.Range("B12").Value = '010268'
$q3.Range("C12").Value = '太平睿安混合A'
$q3.Range("D12").Value = '4.03'
$q3.Range("E12").Value = '39.63'
$q3.Range("F12").Value = '1.61'
$q3.Range("G12").Value = '0.0649'
$q3.Range("H12").Value = 5
$q3.Range("A13").Value = 11
$q3.Range("B13").Value = '159758'
$q3.Range("C13").Value = '华夏中证红利质量ETF'
$q3.Range("D13").Value = '1.69'
$q3.Range("E13").Value = '99.23'
$q3.Range("F13").Value = '3.70'
$q3.Range("G13").Value = '0.0625'
$q3.Range("H13").Value = 5
$q3.Range("A14").Value = 12
$q3.Range("B14").Value = '970023'
$q3.Range("C14").Value = '天风天盈一年定期开放混合'
$q3.Range("D14").Value = '1.24'
$q3.Range("E14").Value = '70.03'
$q3.Range("F14").Value = '4.81'
$q3.Range("G14").Value = '0.0596'
$q3.Range("H14").Value = 4
$q3.Range("A15").Value = 13
$q3.Range("B15").Value = '010060'
$q3.Range("C15").Value = '华泰柏瑞景利混合A'
$q3.Range("D15").Value = '7.16'
$q3.Range("E15").Value = '20.94'
$q3.Range("F15").Value = '0.78'
$q3.Range("G15").Value = '0.0558'
$q3.Range("H15").Value = 9
$q3.Range("A16").Value = 14
$q3.Range("B16").Value = '011383'
$q3.Range("C16").Value = '富安达医药创新混合'
$q3.Range("D16").Value = '1.43'
$q3.Range("E16").Value = '83.09'
$q3.Range("F16").Value = '3.55'
$q3.Range("G16").Value = '0.0508'
$q3.Range("H16").Value = 8
$q3.Range("A17").Value = 15
$q3.Range("B17").Value = '460009'
$q3.Range("C17").Value = '华泰柏瑞量化先行混合A'
$q3.Range("D17").Value = '4.22'
$q3.Range("E17").Value = '93.06'
$q3.Range("F17").Value = '0.97'
$q3.Range("G17").Value = '0.0409'
$q3.Range("H17").Value = 10
$q3.Range("A18").Value = 16
$q3.Range("B18").Value = '005633'
$q3.Range("C18").Value = '建信中证500指数增强C'
$q3.Range("D18").Value = '3.42'
$q3.Range("E18").Value = '82.53'
$q3.Range("F18").Value = '1.12'
$q3.Range("G18").Value = '0.0383'
$q3.Range("H18").Value = 3
$q3.Range("A19").Value = 17
$q3.Range("B19").Value = '013941'
$q3.Range("C19").Value = '东吴医疗服务股票C'
$q3.Range("D19").Value = '0.51'
$q3.Range("E19").Value = '91.62'
$q3.Range("F19").Value = '6.55'
$q3.Range("G19").Value = '0.0334'
$q3.Range("H19").Value = 2
$q3.Range("A20").Value = 18
$q3.Range("B20").Value = '003284'
$q3.Range("C20").Value = '中邮医药健康灵活配置混合'
$q3.Range("D20").Value = '0.65'
$q3.Range("E20").Value = '76.19'
$q3.Range("F20").Value = '3.15'
$q3.Range("G20").Value = '0.0205'
$q3.Range("H20").Value = 8
$q3.Range("A21").Value = 19
$q3.Range("B21").Value = '014344'
$q3.Range("C21").Value = '鹏华中证500指数增强A'
$q3.Range("D21").Value = '1.14'
$q3.Range("E21").Value = '92.67'
$q3.Range("F21").Value = '1.64'
$q3.Range("G21").Value = '0.0187'
$q3.Range("H21").Value = 10
$q3.Range("A22").Value = 20
$q3.Range("B22").Value = '001861'
$q3.Range("C22").Value = '富安达健康人生灵活配置混合A'
$q3.Range("D22").Value = '0.51'
$q3.Range("E22").Value = '86.44'
$q3.Range("F22").Value = '3.62'
$q3.Range("G22").Value = '0.0185'
$q3.Range("H22").Value = 9
$q3.Range("A23").Value = 21
$q3.Range("B23").Value = '014345'
$q3.Range("C23").Value = '鹏华中证500指数增强C'
$q3.Range("D23").Value = '0.73'
$q3.Range("E23").Value = '92.67'
$q3.Range("F23").Value = '1.64'
$q3.Range("G23").Value = '0.0120'
$q3.Range("H23").Value = 10
$q3.Range("A24").Value = 22
$q3.Range("B24").Value = '011824'
$q3.Range("C24").Value = '浙商汇金量化臻选股票A'
$q3.Range("D24").Value = '0.88'
$q3.Range("E24").Value = '92.26'
$q3.Range("F24").Value = '1.22'
$q3.Range("G24").Value = '0.0107'
$q3.Range("H24").Value = 7
$q3.Range("A25").Value = 23
$q3.Range("B25").Value = '005210'
$q3.Range("C25").Value = '东吴双三角股票C'
$q3.Range("D25").Value = '0.10'
$q3.Range("E25").Value = '91.09'
$q3.Range("F25").Value = '6.16'
$q3.Range("G25").Value = '0.0062'
$q3.Range("H25").Value = 2
$q3.Range("A26").Value = 24
$q3.Range("B26").Value = '003242'
$q3.Range("C26").Value = '创金合信量化发现灵活配置混合C'
$q3.Range("D26").Value = '0.40'
$q3.Range("E26").Value = '92.08'
$q3.Range("F26").Value = '1.53'
$q3.Range("G26").Value = '0.0061'
$q3.Range("H26").Value = 7
$q3.Range("A27").Value = 25
$q3.Range("B27").Value = '010061'
$q3.Range("C27").Value = '华泰柏瑞景利混合C'
$q3.Range("D27").Value = '0.70'
$q3.Range("E27").Value = '20.94'
$q3.Range("F27").Value = '0.78'
$q3.Range("G27").Value = '0.0055'
$q3.Range("H27").Value = 9
$q3.Range("A28").Value = 26
$q3.Range("B28").Value = '005209'
$q3.Range("C28").Value = '东吴双三角股票A'
$q3.Range("D28").Value = '0.09'
$q3.Range("E28").Value = '91.09'
$q3.Range("F28").Value = '6.16'
$q3.Range("G28").Value = '0.0055'
$q3.Range("H28").Value = 2
$q3.Range("A29").Value = 27
$q3.Range("B29").Value = '003241'
$q3.Range("C29").Value = '创金合信量化发现灵活配置混合A'
$q3.Range("D29").Value = '0.32'
$q3.Range("E29").Value = '92.08'
$q3.Range("F29").Value = '1.53'
$q3.Range("G29").Value = '0.0049'
$q3.Range("H29").Value = 7
$q3.Range("A30").Value = 28
$q3.Range("B30").Value = '011825'
$q3.Range("C30").Value = '浙商汇金量化臻选股票C'
$q3.Range("D30").Value = '0.39'
$q3.Range("E30").Value = '92.26'
$q3.Range("F30").Value = '1.22'
$q3.Range("G30").Value = '0.0048'
$q3.Range("H30").Value = 7
$q3.Range("A31").Value = 29
$q3.Range("B31").Value = '010269'
$q3.Range("C31").Value = '太平睿安混合C'
$q3.Range("D31").Value = '0.25'
$q3.Range("E31").Value = '39.63'
$q3.Range("F31").Value = '1.61'
$q3.Range("G31").Value = '0.0040'
$q3.Range("H31").Value = 5
$q3.Range("A32").Value = 30
$q3.Range("B32").Value = '010246'
$q3.Range("C32").Value = '华泰柏瑞量化先行混合C'
$q3.Range("D32").Value = '0.25'
$q3.Range("E32").Value = '93.06'
$q3.Range("F32").Value = '0.97'
$q3.Range("G32").Value = '0.0024'
$q3.Range("H32").Value = 10
$q3.Range("A33").Value = 31
$q3.Range("B33").Value = '005966'
$q3.Range("C33").Value = '安信中证500指数增强C'
$q3.Range("D33").Value = '0.16'
$q3.Range("E33").Value = '92.50'
$q3.Range("F33").Value = '1.06'
$q3.Range("G33").Value = '0.0017'
$q3.Range("H33").Value = 8
$q3.Range("A34").Value = 32
$q3.Range("B34").Value = '005965'
$q3.Range("C34").Value = '安信中证500指数增强A'
$q3.Range("D34").Value = '0.10'
$q3.Range("E34").Value = '92.50'
$q3.Range("F34").Value = '1.06'
$q3.Range("G34").Value = '0.0011'
$q3.Range("H34").Value = 8
$q3.Range("A35").Value = 33
$q3.Range("B35").Value = '014328'
$q3.Range("C35").Value = '格林新兴产业混合C'
$q3.Range("D35").Value = '0.03'
$q3.Range("E35").Value = '64.92'
$q3.Range("F35").Value = '3.49'
$q3.Range("G35").Value = '0.0010'
$q3.Range("H35").Value = 9
$q3.Range("A36").Value = 34
$q3.Range("B36").Value = '014470'
$q3.Range("C36").Value = '富安达健康人生灵活配置混合C'
$q3.Range("D36").Value = '0.01'
$q3.Range("E36").Value = '86.44'
$q3.Range("F36").Value = '3.62'
$q3.Range("G36").Value = '0.0004'
$q3.Range("H36").Value = 9
$q3.Range("A37").Value = 35
$q3.Range("B37").Value = '014327'
$q3.Range("C37").Value = '格林新兴产业混合A'
$q3.Range("D37").Value = '0.01'
$q3.Range("E37").Value = '64.92'
$q3.Range("F37").Value = '3.49'
$q3.Range("G37").Value = '0.0003'
$q3.Range("H37").Value = 9

Write-Host "Edit complete"
